$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates ---
# Several of these price strings are purely numeric-looking (e.g. "1.000",
# "2.180", "263.31"). A plain Range.Value assignment would let Excel
# auto-coerce them into actual numbers and silently drop the significant
# trailing zeros / formatting, so each target cell is temporarily switched
# to Text format before the value is written, which keeps it a literal string.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.513.50'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.846.72'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '263.31'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5223'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3230'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06771'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7703'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07776'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.866.40'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '88.35'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.008'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.93'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007936'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '26.567.10'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.091.41'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.611'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.426'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.962'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.20'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.182'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.676'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '111.69'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.161'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08737'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.107'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04817'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.867'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7143'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.101'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01783'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.180'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.4838'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '112.39'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.046'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.000'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.604'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4159'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05900'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.066'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.1224'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.8839'

# The Text-format switch above changes each cells style id; restore the
# original (default/general) look by pasting the formatting only from an
# unrelated, unchanged column-D cell that still carries the default style.
$ws.Range("D42").Copy()
$ws.Range("D2,D3,D5,D7,D8,D9,D11,D12,D13,D14,D15,D17,D19,D20,D21,D22,D23,D24,D25,D26,D27,D29,D30,D31,D32,D33,D35,D36,D37,D38,D39,D40,D41,D43,D44,D45,D46,D47,D48,D50,D51").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Volume(1h) (column E) updates ---
# These are always padded percentage strings, never ambiguous with numbers,
# so a direct value assignment is sufficient.
$ws.Range("E2").Value = '  +0.29%  '
$ws.Range("E3").Value = '  +0.18%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("E5").Value = '  +1.11%  '
$ws.Range("E6").Value = '  +0.08%  '
$ws.Range("E7").Value = '  +1.69%  '
$ws.Range("E8").Value = '  +0.69%  '
$ws.Range("E9").Value = '  +0.22%  '
$ws.Range("E10").Value = '  -0.45%  '
$ws.Range("E11").Value = '  +0.27%  '
$ws.Range("E12").Value = '  +1.20%  '
$ws.Range("E13").Value = '  +1.06%  '
$ws.Range("E14").Value = '  -0.46%  '
$ws.Range("E15").Value = '  -0.33%  '
$ws.Range("E16").Value = '  +0.11%  '
$ws.Range("E17").Value = '  -0.97%  '
$ws.Range("E18").Value = '  +0.06%  '
$ws.Range("E19").Value = '  +0.47%  '
$ws.Range("E20").Value = '  +0.39%  '
$ws.Range("E21").Value = '  -0.43%  '
$ws.Range("E22").Value = '  +0.88%  '
$ws.Range("E23").Value = '  -1.33%  '
$ws.Range("E24").Value = '  +0.13%  '
$ws.Range("E25").Value = '  -1.20%  '
$ws.Range("E26").Value = '  -6.09%  '
$ws.Range("E27").Value = '  +0.98%  '
$ws.Range("E28").Value = '  +0.21%  '
$ws.Range("E29").Value = '  +0.60%  '
$ws.Range("E30").Value = '  -0.34%  '
$ws.Range("E31").Value = '  +0.13%  '
$ws.Range("E32").Value = '  -1.31%  '
$ws.Range("E33").Value = '  +0.05%  '
$ws.Range("E34").Value = '  -0.51%  '
$ws.Range("E35").Value = '  +0.97%  '
$ws.Range("E36").Value = '  +4.10%  '
$ws.Range("E37").Value = '  +0.69%  '
$ws.Range("E38").Value = '  -1.32%  '
$ws.Range("E39").Value = '  -1.04%  '
$ws.Range("E40").Value = '  -1.50%  '
$ws.Range("E41").Value = '  -0.80%  '
$ws.Range("E43").Value = '  -1.52%  '
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("E45").Value = '  -2.06%  '
$ws.Range("E46").Value = '  -1.95%  '
$ws.Range("E47").Value = '  +0.24%  '
$ws.Range("E48").Value = '  -0.43%  '
$ws.Range("E49").Value = '  -0.19%  '
$ws.Range("E50").Value = '  -3.59%  '
$ws.Range("E51").Value = '  +3.67%  '
